# fall 22 week 14 complete
# Append 24 new matchup rows (1479-1502) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,2,4,0),
    @(2,2,3,1),
    @(5,2,5,1),
    @(6,0,3,2),
    @(3,1,7,2),
    @(4,0,4,2),
    @(5,2,7,1),
    @(5,0,6,2),
    @(2,3,3,0),
    @(3,2,3,1),
    @(5,0,4,3),
    @(5,2,6,0),
    @(6,2,4,1),
    @(4,2,6,1),
    @(4,0,3,3),
    @(4,1,4,2),
    @(4,2,5,0),
    @(6,2,6,0),
    @(5,3,5,0),
    @(3,0,5,3),
    @(3,1,5,2),
    @(6,0,7,3),
    @(3,3,2,0),
    @(5,2,4,1)
)

$startRow = 1479
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the view so the new end-of-data cell is selected, matching a
# user who just finished typing at the bottom of the sheet.
$excel.ActiveWindow.ScrollRow = 1476
[void]$ws.Range("A1503").Select()
